$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45872
$ws.Range("B2").Value = 80.40000000000001
$ws.Range("C2").Value = 69.26000000000001
$ws.Range("D2").Value = 42.96
$ws.Range("E2").Value = 37.38
$ws.Range("F2").Value = 30.02
$ws.Range("G2").Value = 30.2
$ws.Range("H2").Value = 32.03
$ws.Range("I2").Value = 26.28
$ws.Range("J2").Value = 6.16
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = -1.01
$ws.Range("M2").Value = -2.1
$ws.Range("N2").Value = -2
$ws.Range("Q2").Value = -0.01
$ws.Range("R2").Value = -1.2
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0.05
$ws.Range("U2").Value = 27.08
$ws.Range("V2").Value = 78.2
$ws.Range("W2").Value = 96.13
$ws.Range("X2").Value = 97.40000000000001
$ws.Range("Y2").Value = 85.09
$ws.Range("Z2").Value = 30.51
$ws.Range("AB2").Value = 89.2
$ws.Range("AD2").Value = 91.23999999999999
$ws.Range("AE2").Value = "20h-22h"
$ws.Range("AF2").Value = 87.16
$ws.Range("AG2").Value = "4h-19h"
